# issue #5: stock data from json to db
# Adds 'category' column (value "normal") right after property_category,
# and appends 'source_file' (value "tmpc08e1") + 'index' (row id) columns
# at the end of the 股票 (stock) worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Insert a new column before the existing "date" column (column I),
# shifting date/legislator_name/legislator_id one column to the right.
$ws.Range("I1").EntireColumn.Insert()

# New column I: category header + "normal" value on every data row.
$ws.Range("I1").Value = "category"
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 9).Value = "normal"
}

# New trailing columns M (source_file) and N (index).
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 13).Value = "tmpc08e1"
    $ws.Cells.Item($r, 14).Value = $ws.Cells.Item($r, 1).Value()
}

# Match the header styling (bold, centered, bordered) used by the other
# header cells on row 1.
$ws.Range("M1:N1").Font.Bold = $true
$ws.Range("M1:N1").HorizontalAlignment = -4108
$ws.Range("M1:N1").VerticalAlignment = -4160
$ws.Range("M1:N1").Borders.LineStyle = 1
